$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential notice date (2021-05-25 -> 2021-05-26)
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-26 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-10
$ws.Range("D2").Value = 0.09238292390281035
$ws.Range("E2").Value = 0.01698540513336688

$ws.Range("D3").Value = 0.1073841722742427
$ws.Range("E3").Value = 0.007432070872227747

$ws.Range("D4").Value = 0.1198814745362923
$ws.Range("E4").Value = 0.00186814042188832

$ws.Range("D5").Value = 0.1403336146800688
$ws.Range("E5").Value = 0.00547358400761544

$ws.Range("D6").Value = 0.1363238629927783
$ws.Range("E6").Value = -0.0006905123601712582

$ws.Range("D7").Value = 0.1456987223555632
$ws.Range("E7").Value = 0.00528551203397809

$ws.Range("D8").Value = 0.1279836515156618
$ws.Range("E8").Value = 0.005589879376286966

$ws.Range("D9").Value = 0.1300115777425827
$ws.Range("E9").Value = 0.01060924580054867

$ws.Range("E10").Value = 0.006130028423512623
